$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6957.25
$ws.Range("I43").Value = 6776.3335
$ws.Range("K43").Value = 6776.3335
$ws.Range("M43").Value = -6707.3335
$ws.Range("H86").Value = 3949
$ws.Range("I86").Value = 3949
$ws.Range("K86").Value = 3949
$ws.Range("M86").Value = -2826
$ws.Range("H89").Value = 3949
$ws.Range("I89").Value = 3949
$ws.Range("K89").Value = 19745
$ws.Range("M89").Value = -14129
$ws.Range("H92").Value = 1102.6666
$ws.Range("I92").Value = 1243.2
$ws.Range("J92").Value = 400
$ws.Range("K92").Value = 1243.2
$ws.Range("L92").Value = 400
$ws.Range("M92").Value = 4.799999999999955
$ws.Range("N92").Value = -2896
$ws.Range("H106").Value = 7850.6665
$ws.Range("I106").Value = 7776
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 7776
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -7145
$ws.Range("N106").Value = -9262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1174.6428
$ws.Range("I2").Value = 1160.091
$ws.Range("J2").Value = 1228
$ws.Range("K2").Value = 1160.091
$ws.Range("L2").Value = 1228
$ws.Range("M2").Value = -1047.091
$ws.Range("N2").Value = -1454
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H62").Value = 60000
$ws.Range("J62").Value = 60000
$ws.Range("L62").Value = 60000
$ws.Range("N62").Value = -61248
$ws.Range("H65").Value = 60000
$ws.Range("J65").Value = 60000
$ws.Range("L65").Value = 180000
$ws.Range("N65").Value = -186240
$ws.Range("H74").Value = 1383.3334
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126
$ws.Range("H77").Value = 1383.3334
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632
$ws.Range("H110").Value = 765.3333
$ws.Range("I110").Value = 679
$ws.Range("K110").Value = 679
$ws.Range("M110").Value = 1366
$ws.Range("H116").Value = 1174.6428
$ws.Range("I116").Value = 1160.091
$ws.Range("J116").Value = 1228
$ws.Range("K116").Value = 1160.091
$ws.Range("L116").Value = 1228
$ws.Range("M116").Value = 1133.909
$ws.Range("N116").Value = -5816
$ws.Range("H132").Value = 3056.4211
$ws.Range("I132").Value = 2007.5454
$ws.Range("K132").Value = 6022.6362
$ws.Range("M132").Value = -3492.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1174.6428
$ws.Range("I3").Value = 1160.091
$ws.Range("J3").Value = 1228
$ws.Range("K3").Value = 1160.091
$ws.Range("L3").Value = 1228
$ws.Range("M3").Value = -1046.091
$ws.Range("N3").Value = -1456
$ws.Range("H94").Value = 1052.7693
$ws.Range("I94").Value = 638.1818
$ws.Range("K94").Value = 638.1818
$ws.Range("M94").Value = -187.1818
$ws.Range("H134").Value = 2457.8125
$ws.Range("I134").Value = 2105
$ws.Range("J134").Value = 2911.4285
$ws.Range("K134").Value = 6315
$ws.Range("L134").Value = 8734.2855
$ws.Range("M134").Value = -3780
$ws.Range("N134").Value = -13804.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1581.625
$ws.Range("I31").Value = 1435.25
$ws.Range("K31").Value = 1435.25
$ws.Range("M31").Value = -1140.25
$ws.Range("H34").Value = 1581.625
$ws.Range("I34").Value = 1435.25
$ws.Range("K34").Value = 1435.25
$ws.Range("M34").Value = -1233.25
$ws.Range("H43").Value = 18000
$ws.Range("J43").Value = 18000
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18368
$ws.Range("H94").Value = 2000
$ws.Range("J94").Value = 2000
$ws.Range("L94").Value = 2000
$ws.Range("N94").Value = -2902
$ws.Range("H99").Value = 4630.4287
$ws.Range("I99").Value = 4566.5
$ws.Range("K99").Value = 4566.5
$ws.Range("M99").Value = -3068.5
$ws.Range("H101").Value = 18000
$ws.Range("J101").Value = 18000
$ws.Range("L101").Value = 18000
$ws.Range("N101").Value = -24490
$ws.Range("H126").Value = 4630.4287
$ws.Range("I126").Value = 4566.5
$ws.Range("K126").Value = 13699.5
$ws.Range("M126").Value = -11229.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H107").Value = 1899.6666
$ws.Range("J107").Value = 279.6
$ws.Range("L107").Value = 838.8000000000001
$ws.Range("N107").Value = -4678.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 90.36364
$ws.Range("J2").Value = 95
$ws.Range("L2").Value = 95
$ws.Range("N2").Value = -321
$ws.Range("H132").Value = 3820.5715
$ws.Range("J132").Value = 3824.4167
$ws.Range("L132").Value = 11473.2501
$ws.Range("N132").Value = -16533.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H46").Value = 2997.25
$ws.Range("I46").Value = 2994.5
$ws.Range("K46").Value = 2994.5
$ws.Range("M46").Value = -2806.5
$ws.Range("H136").Value = 3723.0908
$ws.Range("J136").Value = 2222
$ws.Range("L136").Value = 6666
$ws.Range("N136").Value = -11766

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 197.25
$ws.Range("I113").Value = 144
$ws.Range("K113").Value = 432
$ws.Range("M113").Value = 1738
$ws.Range("H132").Value = 1805.6364
$ws.Range("I132").Value = 1543.7222
$ws.Range("K132").Value = 4631.1666
$ws.Range("M132").Value = -2101.1666
$ws.Range("H136").Value = 2200
$ws.Range("I136").Value = 1760.2174
$ws.Range("J136").Value = 3323.889
$ws.Range("K136").Value = 5280.6522
$ws.Range("L136").Value = 9971.667000000001
$ws.Range("M136").Value = -2730.6522
$ws.Range("N136").Value = -15071.667
